$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Add the new "batch" worksheet right after "programData"
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "batch"

# --- Row 1: headers ---
$ws2.Range("A1").Value = "scenario"
$ws2.Range("B1").Value = "batch-suffix"
$ws2.Range("C1").Value = "status"
$ws2.Range("D1").Value = "no-of-classes"
$ws2.Range("E1").Value = "error-message"

# --- Row 2 ---
$ws2.Range("A2").Value = "VALID"
$ws2.Range("B2").Value = 179
$ws2.Range("B2").Font.Color = 255
$ws2.Range("C2").Value = "ACTIVE"
$ws2.Range("D2").Value = 10

# --- Row 3 ---
$ws2.Range("A3").Value = "INVALID-BATCH-SUFFIX"
$ws2.Range("B3").Value = "asd"
$ws2.Range("C3").Value = "ACTIVE"
$ws2.Range("D3").Value = 10
$ws2.Range("E3").Value = "This field accept only numbers and max 5 count."

# --- Row 4 ---
$ws2.Range("A4").Value = "INVALID-MISSING-MANDATORY-FIELD"
$ws2.Range("B4").Value = 179
$ws2.Range("C4").Value = "ACTIVE"
$ws2.Range("E4").Value = "Number of classes is required."

# --- Row 5 ---
$ws2.Range("A5").Value = "VALID-MANDATORY-FIELDS"
$ws2.Range("B5").Value = 183
$ws2.Range("C5").Value = "ACTIVE"
$ws2.Range("D5").Value = 10

# --- Row 6 ---
$ws2.Range("A6").Value = "VALID-MANDATORY-FIELDS-CANCEL"
$ws2.Range("B6").Value = 210
$ws2.Range("C6").Value = "ACTIVE"
$ws2.Range("D6").Value = 10
$ws2.Range("E6").Value = "Showing 0 to 0 of 0 entries"

# --- Row 7 ---
$ws2.Range("A7").Value = "VALID-EDIT"
$ws2.Range("B7").Value = 183
$ws2.Range("B7").Font.Color = 255
$ws2.Range("C7").Value = "ACTIVE"
$ws2.Range("D7").Value = 20
$ws2.Range("E7").Value = "Successful"

# --- Row 8: a single formatted-but-empty cell ---
$ws2.Range("B8").Font.Name = "Calibri"

# Column sizing / view state to mirror the authored layout (closest achievable
# to the source 19.9140625 "best fit" width given this host's pixel-grid rounding)
$ws2.Columns.Item(1).ColumnWidth = 19

[void]$ws2.Cells.Select()
[void]$ws2.Range("K22").Select()
[void]$ws2.Cells.Select()
[void]$ws2.Activate()
